$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 & 10: turn the old "LOGIN" / "PRODUCTOS" cases into the new
#     "LOGIN (Positivo)" / "HOME (Negativo)" cases, and fill in the now-used
#     "RESULTADO OBTIDO" (D) / "COMENTÁRIOS" (E) columns
$ws.Range("D9").Value = "Usuário ser redirecionado para a Home"
$ws.Range("D10").Value = "Ter apenas alimentos na Home"
$ws.Range("B9").Value = "LOGIN (Positivo)"
$ws.Range("C10").Value = "Cadastrar  um produto diferente de um alimento"
$ws.Range("B10").Value = "HOME (Negativo)"
$ws.Range("C9").Value = "Usuário entrar no Comida Ya! Com credenciais válidas"
$ws.Range("E10").Value = "Foi possível cadastrar produtos diferentes de alimentos."
$ws.Range("E9").Value = "Usuário foi redirecionado para a página Home."

# --- Row 10 number, and rows 11-13 keep their existing text, just renumbered
#     (# column) to stay sequential
$ws.Range("A10").Value = 2
$ws.Range("A11").Value = 3
$ws.Range("A12").Value = 4
$ws.Range("A13").Value = 5

# --- Widen column E so the new, longer "COMENTÁRIOS" text fits
$ws.Columns.Item(5).ColumnWidth = 51

# --- Restore scroll position / selection to match the saved view
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("E9").Select()
